$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.402.64'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("D3").Value = '1.681.37'
$ws.Range("E3").Value = '  +1.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.51'
$ws.Range("E5").Value = '  +4.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5312'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2667'
$ws.Range("E8").Value = '  +3.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06444'
$ws.Range("E9").Value = '  +3.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.24'
$ws.Range("E10").Value = '  +2.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07796'

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.688.68'
$ws.Range("E12").Value = '  +2.32%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.503'
$ws.Range("E13").Value = '  +2.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5642'
$ws.Range("E14").Value = '  +5.03%  '

$ws.Range("D15").Value = '0.0₅8440'
$ws.Range("E15").Value = '  +6.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.06'
$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").Value = '26.442.55'
$ws.Range("E17").Value = '  +1.28%  '

$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.835'
$ws.Range("E18").Value = '  +3.60%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.48'
$ws.Range("E20").Value = '  +3.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.39'
$ws.Range("E21").Value = '  +3.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.397'
$ws.Range("E22").Value = '  +4.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.56'
$ws.Range("E24").Value = '  -3.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1267'
$ws.Range("E25").Value = '  +4.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.494'
$ws.Range("E26").Value = '  +1.94%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.32'
$ws.Range("E27").Value = '  +4.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.434'
$ws.Range("E28").Value = '  +3.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06187'
$ws.Range("E29").Value = '  +3.26%  '

$ws.Range("E30").Value = '  +2.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.551'
$ws.Range("E31").Value = '  +3.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.461'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.710'
$ws.Range("E33").Value = '  +5.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.018'
$ws.Range("E34").Value = '  +4.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.789'
$ws.Range("E35").Value = '  +1.74%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.407'
$ws.Range("E36").Value = '  +1.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5762'
$ws.Range("E37").Value = '  -2.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01646'
$ws.Range("E38").Value = '  +3.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.948'
$ws.Range("E39").Value = '  +1.35%  '

$ws.Range("E40").Value = '  +2.62%  '

$ws.Range("D41").Value = '1.055.82'
$ws.Range("E41").Value = '  -1.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  -0.31%  '

$ws.Range("E43").Value = '  -0.12%  '

$ws.Range("D44").Value = '1.831.89'
$ws.Range("E44").Value = '  +1.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.18'
$ws.Range("E45").Value = '  +4.77%  '

$ws.Range("E46").Value = '  -1.22%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.158'
$ws.Range("E47").Value = '  +1.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  -0.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05201'
$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.062'
$ws.Range("E50").Value = '  +3.32%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4225'
$ws.Range("E51").Value = '  -0.32%  '
